$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
# H1: Experiment_Time -> Avg_Experiment_Time
$ws.Cells.Item(1,8).Value = "Avg_Experiment_Time"

# New header columns I1:O1 (6 new Std_* columns + relocated Obs_Prob)
$newHeaderCols    = @(9, 10, 11, 12, 13, 14, 15)
$newHeaderValues  = @("Std_Total_Rounds", "Std_Expl_Cost", "Std_Expl_Eff", "Std_Round_Time", "Std_Agent_Step_Time", "Std_Experiment_Time", "Obs_Prob")

for ($i = 0; $i -lt $newHeaderCols.Length; $i++) {
  $col = $newHeaderCols[$i]
  $ws.Cells.Item(1, $col).Value = $newHeaderValues[$i]
  # Match the existing bold/centered/bordered header style (s="1")
  $ws.Cells.Item(1, 1).Copy()
  $ws.Cells.Item(1, $col).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# --- Data rows (A2:O13) ---
$row2 = @(1, 1, 54.57, 54.57, 3.128393519999999, 0.16071982, 0.16071982, 8.73086006, 5.502513544792166, 5.502513544792166, 0.316143681776283, 0.01865219949656661, 0.01865219949656661, 1.048595980214898, 0.15)
for ($c = 1; $c -le $row2.Length; $c++) { $ws.Cells.Item(2, $c).Value = $row2[$c-1] }
$row3 = @(1, 1, 87.93, 87.93, 1.95778016, 0.09553492, 0.09553492, 8.245246139999999, 11.95772027387731, 11.95772027387731, 0.2686110744317687, 0.016365110726536, 0.016365110726536, 0.8277653025116386, 0.85)
for ($c = 1; $c -le $row3.Length; $c++) { $ws.Cells.Item(3, $c).Value = $row3[$c-1] }
$row4 = @(2, 1, 32.812, 65.556, 2.63784808, 0.24414562, 0.12207272, 3.9677613, 5.131209860184878, 10.23720982152167, 0.3923558174934527, 0.03910229124459018, 0.01955125644706544, 0.696837848454709, 0.15)
for ($c = 1; $c -le $row4.Length; $c++) { $ws.Cells.Item(4, $c).Value = $row4[$c-1] }
$row5 = @(2, 1, 49.184, 96.1, 1.80626002, 0.1464552, 0.07322764000000001, 3.51306394, 8.537442230532477, 15.84196755052047, 0.2963210268653317, 0.02841550228484797, 0.01420768807481578, 0.4646167245799436, 0.85)
for ($c = 1; $c -le $row5.Length; $c++) { $ws.Cells.Item(5, $c).Value = $row5[$c-1] }
$row6 = @(4, 1, 19.674, 78.604, 2.22525534, 0.3115481999999999, 0.07788686, 1.5121634, 3.764048368477812, 15.03766760151944, 0.4146762848006098, 0.07365618695067118, 0.01841417795721764, 0.3839339673462149, 0.15)
for ($c = 1; $c -le $row6.Length; $c++) { $ws.Cells.Item(6, $c).Value = $row6[$c-1] }
$row7 = @(4, 1, 27.332, 100.198, 1.73004372, 0.19713562, 0.049284, 1.3141698, 5.455024342856989, 16.10831311785955, 0.2783568010981531, 0.04348207049260935, 0.01087051025418645, 0.2553811190379629, 0.85)
for ($c = 1; $c -le $row7.Length; $c++) { $ws.Cells.Item(7, $c).Value = $row7[$c-1] }
$row8 = @(6, 1, 13.56, 81.186, 2.18101796, 0.33046818, 0.05507815999999999, 0.7413027799999999, 2.785303543900256, 16.68965227276986, 0.5109835729405287, 0.09686591800126636, 0.01614413893207427, 0.2534691562569671, 0.15)
for ($c = 1; $c -le $row8.Length; $c++) { $ws.Cells.Item(8, $c).Value = $row8[$c-1] }
$row9 = @(6, 1, 19.47, 98.392, 1.7652824, 0.21811254, 0.036352, 0.6950704000000001, 4.310250168264721, 16.0399781108202, 0.3008167468746387, 0.05256429250204821, 0.008760835901407055, 0.1924025233988903, 0.85)
for ($c = 1; $c -le $row9.Length; $c++) { $ws.Cells.Item(9, $c).Value = $row9[$c-1] }
$row10 = @(8, 1, 9.73, 77.634, 2.306508, 0.32872708, 0.04109106000000001, 0.40209858, 2.206998342402353, 17.59053432376604, 0.6076541539804525, 0.1008925755835848, 0.01261139750266959, 0.1644015861238074, 0.15)
for ($c = 1; $c -le $row10.Length; $c++) { $ws.Cells.Item(10, $c).Value = $row10[$c-1] }
$row11 = @(8, 1, 15.076, 92.942, 1.87937462, 0.2129099, 0.02661376, 0.39581328, 3.795132416678006, 16.88457324385945, 0.3455917036581322, 0.05418020617598779, 0.006772594505436608, 0.1285770873477201, 0.85)
for ($c = 1; $c -le $row11.Length; $c++) { $ws.Cells.Item(11, $c).Value = $row11[$c-1] }
$row12 = @(10, 1, 7.554, 75.152, 2.40879072, 0.31911232, 0.03191128, 0.24486496, 1.922101179930907, 19.1380622544109, 0.6589880415188208, 0.1045250294167226, 0.01045247636140697, 0.1104740221524476, 0.15)
for ($c = 1; $c -le $row12.Length; $c++) { $ws.Cells.Item(12, $c).Value = $row12[$c-1] }
$row13 = @(10, 1, 12.088, 82.482, 2.12623898, 0.18772058, 0.01877218, 0.22441674, 3.563376390813243, 16.42660472177287, 0.4029035554545969, 0.05156518642617257, 0.005156513114751597, 0.0877269887108234, 0.85)
for ($c = 1; $c -le $row13.Length; $c++) { $ws.Cells.Item(13, $c).Value = $row13[$c-1] }
